$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A never-written cell used purely as a "blank" source for copy operations,
# so that destination cells end up as bare empty cells (<c r="..."/>)
# instead of being removed entirely (which is what happens when a plain
# empty string is assigned through .Value).
$blank = $ws.Range("Z100")

# ---- Cells that held data in the old 4-column layout but must become
# ---- empty now that the sheet has been restructured to 9 columns.
$ws.Range("B4").ClearContents()
$ws.Range("D4").ClearContents()
$ws.Range("B6").ClearContents()

# ---- New header block (rows 1-2, columns B/C): written column-by-column
# ---- (B1,B2 then C1,C2) so the shared-string table gets
# ---- "sourcelanguage","en","language","it" in that exact order.
$ws.Range("B1").Value = "sourcelanguage"
$ws.Range("B2").Value = "en"
$ws.Range("C1").Value = "language"
$ws.Range("C2").Value = "it"

# A1 ("TsVersion") and A2 ("2.1") are unchanged from the original file, so
# they are intentionally left untouched (re-assigning "2.1" as a string
# would otherwise be auto-coerced into a numeric value by Excel).

# ---- Header row (row 3) ----
$ws.Range("A3").Value = "Context"
$ws.Range("B3").Value = "ID"
$ws.Range("C3").Value = "Source"
$ws.Range("D3").Value = "Translation"
$ws.Range("E3").Value = "TranslationType"
$ws.Range("F3").Value = "comment"
$ws.Range("G3").Value = "extracomment"
$ws.Range("H3").Value = "translatorcomment"
$ws.Range("I3").Value = "Location"

# ---- Row 4 ----
$ws.Range("A4").Value = "MenuBar"
$blank.Copy($ws.Range("B4"))
$ws.Range("C4").Value = "text"
$blank.Copy($ws.Range("D4"))
$ws.Range("E4").Value = "unfinished"
$ws.Range("F4").Value = "my comment"
$blank.Copy($ws.Range("G4"))
$blank.Copy($ws.Range("H4"))
$ws.Range("I4").Value = "../src/app/qml/MenuBar.qml - 17"

# ---- Row 5 ----
$ws.Range("A5").Value = "MenuBar"
$ws.Range("B5").Value = "1abc"
$ws.Range("C5").Value = "map"
$ws.Range("D5").Value = "whatever"
$ws.Range("E5").Value = "vanished"
$blank.Copy($ws.Range("F5"))
$ws.Range("G5").Value = "extra comment test"
$ws.Range("H5").Value = "translatorcomment"
$ws.Range("I5").Value = "../src/app/qml/MenuBar.qml - 28"

# ---- Row 6 ----
$ws.Range("A6").Value = "MenuBar"
$blank.Copy($ws.Range("B6"))
$ws.Range("C6").Value = "cam"
$ws.Range("D6").Value = "whichever"
$ws.Range("E6").Value = "obsolete"
$blank.Copy($ws.Range("F6"))
$blank.Copy($ws.Range("G6"))
$blank.Copy($ws.Range("H6"))
$ws.Range("I6").Value = "../src/app/qml/MenuBar.qml - 43"

# ---- Row 7 (new row) ----
$ws.Range("A7").Value = "MenuBar"
$blank.Copy($ws.Range("B7"))
$ws.Range("C7").Value = "checklist"
$ws.Range("D7").Value = "mytranslation"
$blank.Copy($ws.Range("E7"))
$blank.Copy($ws.Range("F7"))
$blank.Copy($ws.Range("G7"))
$blank.Copy($ws.Range("H7"))
$ws.Range("I7").Value = "../src/app/qml/MenuBar.qml - 58"
